# DemoSiteDDT.xlsx: the "Result" column (E) for each data-driven test row
# previously held the literal text "true". Update it to read "PASS" instead,
# and register "PASS" as a new shared string used by the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E5").Value = "PASS"
